$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column before column B (existing B->C, C->D) ------------
$ws.Columns("B").Insert() | Out-Null

# Column B should be the same width as column A (75.81640625 characters).
# (ColumnWidth rounds to the engine's internal granularity, which is the
# closest achievable approximation of the source width.)
$ws.Range("B1").ColumnWidth = $ws.Range("A1").ColumnWidth

# --- Row 2 query text (write A2 first so new shared-strings land in the
#     same order as the source commit: bam-query, StatQuery, stat-query) ---
$bamQuery = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report)OPTIONAL MATCH (s)<-[*]-(f:file)  WITH  c AS c, t ,a, s WHERE f.file_format IN ['bam']  RETURN DISTINCT coalesce(c.case_id,'') AS ``Case ID`` , coalesce(t.clinical_trial_designation ,'')as ``Trial Code`` , coalesce(a.arm_id,'') As ``Arm`` , coalesce(a.arm_drug,'') As ``Arm Treatment`` , coalesce(c.disease,'') As Diagnosis , coalesce(c.gender,'') As Gender , coalesce(c.race,'') As Race , coalesce(c.ethnicity,'') As Ethnicity"
$statQuery = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report)OPTIONAL MATCH (s)<-[*]-(f:file)  WITH  c AS c, t ,a, s , f WHERE f.file_format IN ['bai','bam','vcf'] RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"

$ws.Range("A2").Value2 = $bamQuery
$ws.Range("B1").Value2 = "StatQuery"
$ws.Range("B2").Value2 = $statQuery

# A2/B2 share the existing wrap-text style used by A2 originally.
$ws.Range("B2").WrapText = $true

# Row 2 grows to fit the wrapped query text (7 lines tall).
$ws.Range("A2:B2").RowHeight = 101.5

# --- View: selection moves to A2 (also clears the old topLeftCell scroll) --
$ws.Range("A2").Select() | Out-Null
